$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "testcases"
$ws2.Range("A3").Value = "Scenario"
$ws2.Range("B3").Value = "Testcase"
$ws2.Range("B15").Select()
$ws1.Select()
